# "works on RW failed script fixations"
# Rows 182-234 held scenarios whose RW ("ready to work"/regression) scripts
# had previously been flagged as failing - column D ("FixationStatus"-style
# Yes/No flag) was "No" for all of them. The fixations now work, so flip
# that flag to "Yes" for every one of those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D182:D234").Value = "Yes"

# Leave the trail showing where the author ended up after scrolling
# through / reviewing the sheet (matches the saved cursor position).
$ws.Activate() | Out-Null
$ws.Range("D237").Select() | Out-Null
